$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.38094952709161589
$ws.Range("A2").Value = -0.009999999637550161
$ws.Range("A3").Value = -0.0089999996318343989
$ws.Range("A4").Value = -0.011999999898019809
$ws.Range("A5").Value = -0.0059999996382789078
$ws.Range("A6").Value = -0.0059999996274342493
$ws.Range("A7").Value = -0.019999999552398506
$ws.Range("A8").Value = -0.019999999547315461
$ws.Range("A9").Value = -0.0059999996142998668
$ws.Range("A10").Value = -0.0059999996086617102
$ws.Range("A11").Value = -0.0044999996156818156
$ws.Range("A12").Value = -0.0059999996063444527
$ws.Range("A13").Value = -0.0059999996003394784
$ws.Range("A14").Value = -0.01199999956709874
$ws.Range("A15").Value = -0.0059999995971162789
$ws.Range("A16").Value = -0.0059999995958497365
$ws.Range("A17").Value = -0.0059999995941266704
$ws.Range("A18").Value = -0.008999999578027662
$ws.Range("A19").Value = -0.0089999996418548278
$ws.Range("A20").Value = -0.0089999996308884889
$ws.Range("A21").Value = -0.067184566937994017
$ws.Range("A22").Value = -0.0089999996280858419
$ws.Range("A23").Value = -0.0089999996285481387
$ws.Range("A24").Value = -0.041999999450112568
$ws.Range("A25").Value = -0.041999999447164704
$ws.Range("A26").Value = -0.021642298145721384
$ws.Range("A27").Value = -0.0059999996267672273
$ws.Range("A28").Value = -0.0059999996264306077
$ws.Range("A29").Value = -0.011999999595223798
$ws.Range("A30").Value = -0.019999999553996783
$ws.Range("A31").Value = -0.014999999583093171
$ws.Range("A32").Value = 0.055691841994184088
$ws.Range("A33").Value = -0.0059999996313200299
